# DynamicRLS.xlsx update:
#  - Add two new users (Joshua Hernandez / Josh BaconRox) to the "Users Table" sheet,
#    each with a mailto hyperlink on their e-mail address.
#  - Add the same two users (by e-mail) to the "UserRegion Table" sheet, mapped to
#    Western Region / Central Region respectively, and mailto-hyperlink several
#    pre-existing e-mail cells that weren't hyperlinks yet.
#  - Resize both Excel Tables to cover the newly added rows.
#  - Leave "UserRegion Table" as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. "Users Table" sheet - append the two new users
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users Table")

$wsUsers.Range("A12").Value = "joshua.hernandez@comparex.com"
$wsUsers.Range("B12").Value = "Joshua Hernandez"
$wsUsers.Range("A13").Value = "student@baconrox.onmicrosoft.com"
$wsUsers.Range("B13").Value = "Josh BaconRox"

$wsUsers.Hyperlinks.Add($wsUsers.Range("A12"), "mailto:joshua.hernandez@comparex.com") | Out-Null
$wsUsers.Hyperlinks.Add($wsUsers.Range("A13"), "mailto:student@baconrox.onmicrosoft.com") | Out-Null

# Match the "UserName" column formatting used by the rest of the table (copy
# from the last pre-existing data row).
$wsUsers.Range("B11").Copy() | Out-Null
$wsUsers.Range("B12").PasteSpecial($xlPasteFormats) | Out-Null
$wsUsers.Range("B13").PasteSpecial($xlPasteFormats) | Out-Null

$loUsers = $wsUsers.ListObjects.Item(1)
$loUsers.Resize($wsUsers.Range("A1:B13")) | Out-Null

# ---------------------------------------------------------------------------
# 2. "UserRegion Table" sheet
# ---------------------------------------------------------------------------
$wsUserRegion = $wb.Worksheets.Item("UserRegion Table")

# Hyperlink a handful of existing e-mail cells that weren't linked yet.
$wsUserRegion.Hyperlinks.Add($wsUserRegion.Range("A4"), "mailto:" + $wsUserRegion.Range("A4").Value2) | Out-Null
$wsUserRegion.Hyperlinks.Add($wsUserRegion.Range("A6"), "mailto:" + $wsUserRegion.Range("A6").Value2) | Out-Null
$wsUserRegion.Hyperlinks.Add($wsUserRegion.Range("A10"), "mailto:" + $wsUserRegion.Range("A10").Value2) | Out-Null
$wsUserRegion.Hyperlinks.Add($wsUserRegion.Range("A14"), "mailto:" + $wsUserRegion.Range("A14").Value2) | Out-Null

# New rows for the two new users.
$wsUserRegion.Range("A17").Value = "joshua.hernandez@comparex.com"
$wsUserRegion.Range("B17").Value = "Western Region"
$wsUserRegion.Range("A18").Value = "student@baconrox.onmicrosoft.com"
$wsUserRegion.Range("B18").Value = "Central Region"

$wsUserRegion.Hyperlinks.Add($wsUserRegion.Range("A17"), "mailto:joshua.hernandez@comparex.com") | Out-Null
$wsUserRegion.Hyperlinks.Add($wsUserRegion.Range("A18"), "mailto:student@baconrox.onmicrosoft.com") | Out-Null

# All hyperlinked "User" cells in this sheet share the same look (font +
# left/top alignment) already used by A2/A3 - make sure every newly
# hyperlinked cell (existing rows and the two brand new ones) matches it.
$wsUserRegion.Range("A2").Copy() | Out-Null
$wsUserRegion.Range("A4").PasteSpecial($xlPasteFormats) | Out-Null
$wsUserRegion.Range("A6").PasteSpecial($xlPasteFormats) | Out-Null
$wsUserRegion.Range("A10").PasteSpecial($xlPasteFormats) | Out-Null
$wsUserRegion.Range("A14").PasteSpecial($xlPasteFormats) | Out-Null
$wsUserRegion.Range("A17").PasteSpecial($xlPasteFormats) | Out-Null
$wsUserRegion.Range("A18").PasteSpecial($xlPasteFormats) | Out-Null

$loUserRegion = $wsUserRegion.ListObjects.Item(1)
$loUserRegion.Resize($wsUserRegion.Range("A1:B18")) | Out-Null

# ---------------------------------------------------------------------------
# 3. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------------
$wsUsers.Range("G19").Select() | Out-Null
$wsUserRegion.Activate() | Out-Null
$wsUserRegion.Range("F14").Select() | Out-Null
